$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix column B (Japanese) cell text where the underlying sentence changed ---
$ws.Cells.Item(86, 2).Value = '発光部305は、複数の発光素子306-1から発光素子306-nを備える(以下、「発光素子306」という。)。'
$ws.Cells.Item(87, 2).Value = '発光素子306は、発光制御部314からの発光指示に応じて発光する。 強制発光させるデータは、発光指示の一態様である。'
$ws.Cells.Item(88, 2).Value = '発光素子306は、画像形成処理中に発光することで、感光体ドラム320に静電潜像を形成させる。'
$ws.Cells.Item(89, 2).Value = '発光素子306は、有機EL(Electro Luminescence)であっても良いし、LED(Light Emitting Diode)であっても良い。'
$ws.Cells.Item(90, 2).Value = '発光素子306は、感光体ドラム320に静電潜像を形成できるならばどのような発光素子であっても良い。 感光体ドラム320は、発光素子306の発光方向に設置される。'
$ws.Cells.Item(91, 2).Value = 'FIG. 4は、実施形態の発光時間記憶部301が記憶する発光素子306の各累積発光時間を記憶する一具体例を示す図である。'
$ws.Cells.Item(92, 2).Value = 'FIG. 4に示される例では、発光素子の値と累積発光時間(分)の値とが同一レコードに記憶される。'
$ws.Cells.Item(93, 2).Value = '発光時間記憶部301の最上段のレコードは、発光素子の値が"306-1"、累積発光時間(分)の値が"300"である。'
$ws.Cells.Item(94, 2).Value = '従って、発光素子306-1は、これまで300分発光したことを示す。'
$ws.Cells.Item(95, 2).Value = 'FIG. 5及びFIG. 6は、実施形態の発光素子306を発光させるか否かを判定するフローチャートである。'
$ws.Cells.Item(96, 2).Value = '制御部310は、強制発光閾値を取得する(ACT101)。'
$ws.Cells.Item(97, 2).Value = '制御部310は、最長記憶部303に0の値を設定する(ACT102)。'
$ws.Cells.Item(98, 2).Value = '制御部310は、最短記憶部302が記録できる最大値(例えば、MAX)を最短記憶部302に設定する(ACT103)。'
$ws.Cells.Item(99, 2).Value = '制御部310は、画像データの主走査方向の画素数をカウントする変数nに1を設定する(ACT104)。'
$ws.Cells.Item(100, 2).Value = '制御部310は、n番目の画素に対応する発光素子306の累積発光時間の値を発光時間記憶部301から取得する(ACT105)。'
$ws.Cells.Item(101, 2).Value = '以下、制御部310が取得した主走査方向からn番目の発光素子306の累積発光時間を「発光時間」という。'
$ws.Cells.Item(102, 2).Value = '制御部310は、主走査方向からn番目の1画素を画像データから取得する(ACT106)。'
$ws.Cells.Item(103, 2).Value = '制御部310は、取得された画像データの1画素が白であるか否か判定する(ACT107)。'
$ws.Cells.Item(104, 2).Value = '取得された画像データの1画素が白でない場合(ACT107: NO)、発光制御部314は、発光時間の値に1を加算する(ACT108)。'
$ws.Cells.Item(105, 2).Value = '取得された画像データの1画素が白である場合(ACT107: YES)、なんらの処理も実行しない。'
$ws.Cells.Item(106, 2).Value = '発光時間判定部311は、発光時間の値が最長記憶部303に記憶される累積発光時間の値よりも大きいか否か判定する(ACT109)。'
$ws.Cells.Item(107, 2).Value = '発光時間の値が最長記憶部303に記憶される累積発光時間の値よりも大きい場合(ACT109: YES)、発光時間判定部311は、最長記憶部303に発光時間の値を記録する(ACT110)。'
$ws.Cells.Item(108, 2).Value = '発光時間の値が最長記憶部303に記憶される累積発光時間の値よりも大きくない場合(ACT109: NO)、なんらの処理も実行しない。'
$ws.Cells.Item(109, 2).Value = '発光時間判定部311は、発光時間と最短記憶部302に記憶する累積発光時間とを比較して最短記憶部302に記憶する累積発光時間の方が小さいかを判定する(ACT111)。'
$ws.Cells.Item(110, 2).Value = '発光時間の値が最短記憶部302に記憶される累積発光時間の値よりも小さい場合(ACT111: YES)、発光時間判定部311は、最短記憶部302に発光時間の値を記録する(ACT112)。'
$ws.Cells.Item(111, 2).Value = '発光素子306の累積発光時間が最短記憶部302に記憶される累積発光時間の値よりも小さくない場合(ACT111: NO)、なんらの処理も実行しない。'
$ws.Cells.Item(112, 2).Value = '発光制御部314は、発光時間の値を発光時間記憶部301の発光素子306-nの累積発光時間の値として記録させる(ACT113)。'
$ws.Cells.Item(113, 2).Value = '制御部310は、nが主走査方向の最終画素まで到達したか否か判定する(ACT114)。'
$ws.Cells.Item(114, 2).Value = '主走査方向の最終画素まで到達していない場合(ACT114: NO)、制御部310は、nに1を加算する(ACT115)。'
$ws.Cells.Item(115, 2).Value = 'ACT115が終了すると、処理はACT105へ遷移する(ACT116)。'
$ws.Cells.Item(116, 2).Value = '主走査方向の最終画素まで到達している場合(ACT114: YES)、制御部310は印刷ジョブが終了したか否か判定する(ACT117)。'
$ws.Cells.Item(117, 2).Value = '印刷ジョブが終了していない場合(ACT117: NO)、制御部310は、発光部305を副走査方向へ移動させる(ACT118)。'
$ws.Cells.Item(118, 2).Value = 'ACT115が終了すると、処理はACT102へ遷移する(ACT119)。'
$ws.Cells.Item(119, 2).Value = '印刷ジョブが終了している場合(ACT117: YES)、差分算出部312は、最長記憶部303に記憶される累積発光時間の値と最短記憶部302に記憶される累積発光時間の値との差分値を算出する(ACT120)。'
$ws.Cells.Item(120, 2).Value = '差分算出部312は、算出された差分値が強制発光閾値よりも小さいか否か判定する(ACT121)。'
$ws.Cells.Item(121, 2).Value = '算出された差分値が強制発光閾値よりも小さい場合(ACT121: YES)、なんらの処理も実行しない。'
$ws.Cells.Item(122, 2).Value = '算出された差分値が強制発光閾値よりも小さくない場合(ACT121: NO)、制御部310は、強制発光フラグを成立させる(ACT122)。'
$ws.Cells.Item(123, 2).Value = 'FIG. 7及びFIG. 8は、実施形態の発光素子306を強制発光させる場合の処理の流れを示すフローチャートである。'
$ws.Cells.Item(124, 2).Value = '制御部310は、調整時間の値を取得する(ACT201)。'
$ws.Cells.Item(125, 2).Value = '制御部310は、調整時間の値を最短記憶部302に設定する(ACT202)。'
$ws.Cells.Item(126, 2).Value = '制御部310は、画像データの主走査方向の画素数をカウントする変数nに1を設定する(ACT203)。'
$ws.Cells.Item(127, 2).Value = '制御部310は、主走査方向からn番目の画素に対応する発光素子306の累積発光時間の値(以下、「発光時間」という。)を発光時間記憶部301から取得する(ACT204)。'
$ws.Cells.Item(128, 2).Value = '発光制御部314は、発光時間の値が調整時間の値よりも大きいか否か判定する(ACT205)。'
$ws.Cells.Item(129, 2).Value = '調整時間の値よりも大きい場合(ACT205: YES)、発光制御部314は、発光素子306を発光させないことを表す非発光のデータを生成する(ACT206)。'
$ws.Cells.Item(130, 2).Value = '調整時間の値よりも大きくない場合(ACT205: NO)、発光制御部314は、発光時間の値に1を加算する(ACT207)。'
$ws.Cells.Item(131, 2).Value = '発光時間判定部311は、最短記憶部302に記憶される累積発光時間が発光時間の値よりも大きいか否か判定する(ACT208)。'
$ws.Cells.Item(132, 2).Value = '発光時間の値よりも大きい場合(ACT208: YES)、発光時間判定部311は、最短記憶部302に発光時間の値を記録させる(ACT209)。'
$ws.Cells.Item(133, 2).Value = '発光素子306の累積発光時間の値よりも大きくない場合(ACT208: YES)、なんらの処理も実行しない。'
$ws.Cells.Item(134, 2).Value = '発光制御部314は、発光時間の値を発光時間記憶部301の発光素子306-nの累積発光時間の値として記録させる(ACT210)。'
$ws.Cells.Item(135, 2).Value = '発光制御部314は、発光素子306を発光させることを表す発光のデータを生成する(ACT211)。'
$ws.Cells.Item(136, 2).Value = '制御部310は、nが主走査方向の最終画素まで到達したか否か判定する(ACT212)。'
$ws.Cells.Item(137, 2).Value = '主走査方向の最終画素まで到達していない場合(ACT212: NO)、制御部310は、nに1を加算する(ACT213)。'
$ws.Cells.Item(138, 2).Value = '主走査方向の最終画素まで到達している場合(ACT212: YES)、発光制御部314は、発光部305に生成されたデータを送信し、発光素子306を発光させる(ACT214)。'
$ws.Cells.Item(139, 2).Value = '発光制御部314は、最短記憶部302に記憶される累積時間の値と調整時間の値が等しいか否か判定する(ACT215)。'
$ws.Cells.Item(140, 2).Value = '等しくない場合(ACT215: NO)、処理はACT202へ遷移する(ACT216)。'
$ws.Cells.Item(141, 2).Value = '等しい場合(ACT216: YES)、調整時間算出部315は、調整時間を算出する(ACT217)。'
$ws.Cells.Item(142, 2).Value = '調整時間算出部315は、調整時間記憶部304に調整時間を記憶させ、処理を終了する(ACT218)。'
$ws.Cells.Item(143, 2).Value = 'FIG. 9は、画像データ(A)を印字した場合の発光素子306の累積発光時間を表した図である。'
$ws.Cells.Item(144, 2).Value = '画像データごとに、発光素子306の累積発光時間が異なる。'
$ws.Cells.Item(145, 2).Value = '画像データ(A)には、主走査方向の1画素目から最終画素まで伸びる黒い画素がある。'
$ws.Cells.Item(146, 2).Value = 'このため、発光部305が備える全ての発光素子306が発光する。'
$ws.Cells.Item(147, 2).Value = 'OLEDヘッドは発光部305の一態様である。 FIG. 10は、画像データ(A)を5部印字後の発光素子306の累積発光時間を示す図である。'
$ws.Cells.Item(148, 2).Value = '累積発光時間の最大値と累積発光時間の最小値の差分値が強制発光閾値よりも大きい場合、調整時間よりも累積発光時間が短い発光素子306は強制発光される。'
$ws.Cells.Item(149, 2).Value = 'FIG. 11は、画像データ(A)を5部印字後に発光素子306の累積発光時間が調整時間に到達した状態を示す図である。'
$ws.Cells.Item(150, 2).Value = '発光制御部314は、調整時間よりも累積発光時間が短い発光素子306は、調整時間まで発光させる。'
$ws.Cells.Item(151, 2).Value = 'これによって、全ての発光素子306が調整時間以上の累積発光時間となる。'
$ws.Cells.Item(152, 2).Value = 'FIG. 12は、画像データ(B)を印字した場合の発光素子306の累積発光時間を表した図である。'
$ws.Cells.Item(153, 2).Value = '画像データ(A)を印字した場合と異なり、画像データ(B)を印字しても、発光しない発光素子306がある。'
$ws.Cells.Item(154, 2).Value = 'FIG. 13は、画像データ(B)を15部印字後の発光素子306の累積発光時間を示す図である。'
$ws.Cells.Item(155, 2).Value = '発光素子306は、画像データ(A)を5部印字した直後の発光素子306が用いられる。 累積発光時間の最大値と累積発光時間の最小値の差分値が強制発光閾値よりも大きい場合、調整時間よりも累積発光時間が短い発光素子306は強制発光される。'
$ws.Cells.Item(156, 2).Value = 'FIG. 14は、画像データ(B)を15部印字後に発光素子306の累積発光時間が調整時間に到達した状態を示す図である。'
$ws.Cells.Item(157, 2).Value = '調整時間よりも累積発光時間が短い発光素子306は、調整時間まで発光する。'
$ws.Cells.Item(158, 2).Value = 'これによって、全ての発光素子306が調整時間以上の累積発光時間となる。'
$ws.Cells.Item(159, 2).Value = 'このように、発光制御部314は発光時間記憶部301に記憶される累積発光時間の値が調整時間の値よりも大きいか否か判定する。'
$ws.Cells.Item(160, 2).Value = '判定の結果、調整時間の値よりも累積発光時間の値が大きくない場合、発光制御部314は発光素子306を発光させるように制御する。'
$ws.Cells.Item(161, 2).Value = 'これに対して、調整時間の値よりも累積発光時間の値が大きい場合、発光制御部314は発光素子306を発光させないように制御する。'
$ws.Cells.Item(162, 2).Value = 'これによって、発光制御部314は、調整時間の値よりも累積発光時間の値が大きくない発光素子306を調整時間の値まで発光させる。'
$ws.Cells.Item(163, 2).Value = 'したがって、各発光素子306を累積発光時間の値が最も大きい累積発光時間の値まで発光させるよりも、累積発光時間の値と調整時間の値との差分値ほど発光時間が短くなる。'
$ws.Cells.Item(164, 2).Value = 'そのため、発光素子306は、発光量を均一に維持しつつ、差分値だけ発光による劣化が抑制される。'
$ws.Cells.Item(165, 2).Value = 'While certain embodiments have been described these embodiments have been presented by way of example only, and are not intended to limit the scope of the inventions.Indeed, the novel embodiments described herein may be embodied in a variety of other forms: furthermore various omissions, substitutions and changes in the form of the embodiments described herein may be made without departing from the spirit of the inventions.The accompanying claims and their equivalents are intended to cover such forms or modifications as would fall within the scope and spirit of the invention. WHAT IS CLAIMED IS: 1. 感光体ドラムに静電潜像を形成させる複数の発光素子を備える発光部と、 累積発光時間が最長である発光素子の累積発光時間よりも短い時間である調整時間を記憶する調整時間記憶部と、 前記複数の発光素子のうち、一の発光素子の累積発光時間と他の発光素子の累積発光時間の差が所定の条件を満たす場合、前記調整時間記憶部に記憶する調整時間よりも累積発光時間が短い発光素子に対して前記調整時間まで前記発光素子を発光させるように制御する発光制御部と、 を備える画像形成装置。 2. クレーム1に記載の画像形成装置であって、 前記発光部が備える複数の発光素子毎の累積発光時間を記憶する発光時間記憶部をさらに備える。'
$ws.Cells.Item(166, 2).Value = 'ABSTRACT 実施形態の画像形成装置は、発光部と、調整時間記憶部と、発光制御部とを持つ。'
$ws.Cells.Item(167, 2).Value = '発光部は、感光体ドラムに静電潜像を形成させる複数の発光素子を備える。'
$ws.Cells.Item(168, 2).Value = '調整時間記憶部は、累積発光時間が最長である発光素子の累積発光時間よりも短い時間である調整時間を記憶する。'

# --- Fix column C (English) cell text where the underlying sentence changed ---
$ws.Cells.Item(87, 3).Value = 'Each of the light emitting elements 306 is connected to the light emission control unitThe light emitting elements 306 emit light in response to a light emission instruction from the light emission control unitData which makes light to be forcibly emitted is an aspect of the light emission instruction.'
$ws.Cells.Item(88, 3).Value = 'The light emitting elements 306 forms electrostatic latent image in a photosensitive drum 320 by making light to be emitted during image formation processing.'
$ws.Cells.Item(89, 3).Value = 'The light emitting elements 306 may be organic electro luminescence (EL), or may be a light emitting diode (LED).'
$ws.Cells.Item(90, 3).Value = 'The light emitting elements 306 may be any type of photosensitive element as long as the electrostatic latent image can be formed in photosensitive drumThe photosensitive drum 320 is installed in a light emission direction of the light emitting elementsThe photosensitive drum 320 forms the electrostatic latent image on the basis of the image information.'
$ws.Cells.Item(91, 3).Value = 'FIG.4 is a diagram illustrating one specific example storing each cumulative light emission time of a light emitting element 306 which is stored in a light emission time storage unit 301 according to the embodiment.'
$ws.Cells.Item(92, 3).Value = 'In the example illustrated in FIG.4, a value of the light emitting element and a value of the cumulative light emission time (minutes) are stored in the same record.'
$ws.Cells.Item(93, 3).Value = 'The uppermost record of the light emission time storage unit 301 shows a value of light emitting element of "306-1" and a value of cumulative light emission time (minutes) of "300".'
$ws.Cells.Item(94, 3).Value = 'Hence, the light emitting element 306-1 indicates that light is emitted for 300 minutes until now.'
$ws.Cells.Item(95, 3).Value = 'FIG.5 and FIG.6 are flowcharts which determine whether or not to cause the light emitting element 306 according to the embodiment to emit light.'
$ws.Cells.Item(96, 3).Value = 'The control unit 310 acquires a forced light emission threshold (ACT101).'
$ws.Cells.Item(97, 3).Value = 'The control unit 310 sets a value of zero to the longest storage unit 303 (ACT102).'
$ws.Cells.Item(98, 3).Value = 'The control unit 310 sets a maximum value (for example, MAX) that the shortest storage unit 302 can store to the shortest storage unit 302 (ACT103).'
$ws.Cells.Item(99, 3).Value = 'The control unit 310 sets a variable n which counts the number of pixels in a main scan direction of image data to "1" (ACT104).'
$ws.Cells.Item(100, 3).Value = 'The control unit 310 acquires a value of the cumulative light emission time of the light emitting elements 306 corresponding to an nth pixel from the light emission time storage unit 301 (ACT105).'
$ws.Cells.Item(101, 3).Value = 'Hereinafter, the cumulative light emission time of the nth light emitting elements 306 in the main scan direction that the control unit 310 acquires is referred to as "light emission time".'
$ws.Cells.Item(102, 3).Value = 'The control unit 310 acquires an nth pixel in the main scan direction from image data (ACT106).'
$ws.Cells.Item(103, 3).Value = 'The control unit 310 determines whether or not the acquired one pixel of the image data is white (ACT107).'
$ws.Cells.Item(104, 3).Value = 'In a case where the acquired one pixel of the image data is not white (ACT107: NO), the light emission control unit 314 adds "1" to the value light emission time (ACT108).'
$ws.Cells.Item(105, 3).Value = 'In a case where the acquired one pixel of the image data is white (ACT107: YES), no processing is performed.'
$ws.Cells.Item(106, 3).Value = 'The light emission time determination unit 311 determines whether or not the value of light emission time is larger than the value of cumulative light emission time stored in the longest storage unit 303 (ACT109).'
$ws.Cells.Item(107, 3).Value = 'In a case where the value of light emission time is larger than the value of cumulative light emission time stored in the longest storage unit 303 (ACT109: YES), the light emission time determination unit 311 stores the value of light emission time in the longest storage unit 303 (ACT110).'
$ws.Cells.Item(108, 3).Value = 'In a case where the value of light emission time is not larger than the value of cumulative light emission time stored in the longest storage unit 303 (ACT109: NO), no processing is performed.'
$ws.Cells.Item(109, 3).Value = 'The light emission time determination unit 311 compares the light emission time and the cumulative light emission time stored in the shortest storage unit 302 and determines whether or not the cumulative light emission time stored in the shortest storage unit 302 is shorter (ACT111).'
$ws.Cells.Item(110, 3).Value = 'In a case where a value of the light emission time is smaller than a value of the cumulative light emission time stored in the shortest storage unit 302 (ACT111: YES), the light emission time determination unit 311 stores the value of the light emission time in the shortest storage unit 302 (ACT112).'
$ws.Cells.Item(111, 3).Value = 'In a case where the cumulative light emission time of the light emitting elements 306 is not smaller than the value of the cumulative light emission time stored in the shortest storage unit 302 (ACT111: NO), no processing is performed.'
$ws.Cells.Item(112, 3).Value = 'The light emission control unit 314 stores the value of the light emission time as the value of the cumulative light emission time of the light emitting elements 306-n of the light emission time storage unit 301 (ACT113).'
$ws.Cells.Item(113, 3).Value = 'The control unit 310 determines whether or not n reaches the final pixel in the main scan direction (ACT114).'
$ws.Cells.Item(114, 3).Value = 'In a case where n does not reach the final pixel in the main scan direction (ACT114: NO), the control unit 310 adds "1" to n (ACT115).'
$ws.Cells.Item(115, 3).Value = 'If ACT115 ends, processing moves to ACT105 (ACT116).'
$ws.Cells.Item(116, 3).Value = 'In a case where n reaches the final pixel in the main scan direction (ACT114: YES), the control unit 310 determines whether or not print job ends (ACT117).'
$ws.Cells.Item(117, 3).Value = 'In a case where the print job does not end (ACT117: NO), the control unit 310 moves the light emitting unit 305 in a sub scan direction (ACT118).'
$ws.Cells.Item(118, 3).Value = 'If ACT115 ends, processing moves to ACT102 (ACT119).'
$ws.Cells.Item(119, 3).Value = 'If the print job ends (ACT117: YES), the difference calculation unit 312 calculates a difference value between a value of the cumulative light emission time stored in the longest storage unit 303 and a value of the cumulative light emission time stored in the shortest storage unit 302 (ACT120).'
$ws.Cells.Item(120, 3).Value = 'The difference calculation unit 312 determines whether or not the calculated difference value is smaller than the forced light emission threshold (ACT121).'
$ws.Cells.Item(121, 3).Value = 'In a case where the calculated difference value is smaller than the forced light emission threshold (ACT121: YES), no processing is performed.'
$ws.Cells.Item(122, 3).Value = 'In a case where the calculated difference value is not smaller than the forced light emission threshold (ACT121: NO), the control unit 310 establishes a forced light emission flag (ACT122).'
$ws.Cells.Item(123, 3).Value = 'FIG.7 and FIG.8 are flowcharts illustrating a flow of processing in a case where the light emitting element 306 according to the embodiment forcibly emits light.'
$ws.Cells.Item(124, 3).Value = 'The control unit 310 acquires a value of the adjustment time (ACT201).'
$ws.Cells.Item(125, 3).Value = 'The control unit 310 sets a value of adjustment time to the shortest storage unit 302 (ACT202).'
$ws.Cells.Item(126, 3).Value = 'The control unit 310 sets "1" as the variable n which counts the number of pixels of the image data in the main scan direction (ACT203).'
$ws.Cells.Item(127, 3).Value = 'The control unit 310 acquires a value of the cumulative light emission time (hereinafter, referred to as "light emission time") of the light emitting elements 306 corresponding to an nth pixel in main scan direction from the light emission time storage unit 301 (ACT204).'
$ws.Cells.Item(128, 3).Value = 'The light emission control unit 314 determines whether or not a value of the light emission time is larger than a value of the adjustment time (ACT205).'
$ws.Cells.Item(129, 3).Value = 'In a case where the value of the light emission time is larger than the value of the adjustment time (ACT205: YES), the light emission control unit 314 generates non-light emission data which indicates that the light emitting elements 306 does not emit light (ACT206).'
$ws.Cells.Item(130, 3).Value = 'In a case where the value of the light emission time is not larger than the value of the adjustment time (ACT205: NO), the light emission control unit 314 adds "1" to the value of the light emission time (ACT207).'
$ws.Cells.Item(131, 3).Value = 'The light emission time determination unit 311 determines whether or not the cumulative light emission time stored in the shortest storage unit 302 is larger than the value of the light emission time (ACT208).'
$ws.Cells.Item(132, 3).Value = 'In a case where the cumulative light emission time is larger than the value of the light emission time (ACT208: YES), the light emission time determination unit 311 stores the value of the light emission time in the shortest storage unit 302 (ACT209).'
$ws.Cells.Item(133, 3).Value = 'In a case where the cumulative light emission time is not larger than the value of the light emission time (ACT208: NO), no processing is performed.'
$ws.Cells.Item(134, 3).Value = 'The light emission control unit 314 stores the value of the light emission time as a value of cumulative light emission time of the light emitting element 306-n of the light emission time storage unit 301 (ACT210).'
$ws.Cells.Item(135, 3).Value = 'The light emission control unit 314 generates light emission data indicating that the light emitting elements 306 emits light (ACT211).'
$ws.Cells.Item(136, 3).Value = 'The control unit 310 determines whether or not n reaches the final pixel in the main scan direction (ACT212).'
$ws.Cells.Item(137, 3).Value = 'In a case where n does not reach the final pixel in the main scan direction (ACT212: NO), the control unit 310 adds "1" to n (ACT213).'
$ws.Cells.Item(138, 3).Value = 'In a case where n reaches the final pixel in the main scan direction (ACT212: YES), the light emission control unit 314 transmits data which is generated to the light emitting unit 305 and makes the light emitting elements 306 emit light (ACT214).'
$ws.Cells.Item(139, 3).Value = 'The light emission control unit 314 determines whether or not the value of the cumulative light emission time stored in the shortest storage unit 302 is equal to the value of the adjustment time (ACT215).'
$ws.Cells.Item(140, 3).Value = 'In a case where they are not equal (ACT215: NO), processing moves to ACT202 (ACT216).'
$ws.Cells.Item(141, 3).Value = 'In a case where they are equal (ACT216: YES), the adjustment time calculation unit 315 calculates adjustment time (ACT217).'
$ws.Cells.Item(142, 3).Value = 'The adjustment time calculation unit 315 stores the adjustment time in the adjustment time storage unit 304 and ends processing (ACT218).'
$ws.Cells.Item(143, 3).Value = 'FIG.9 is a diagram illustrating the cumulative light emission time of the light emitting element 306 in a case where image data (A) is printed.'
$ws.Cells.Item(144, 3).Value = 'The cumulative light emission time of the light emitting elements 306 differs in each image data.'
$ws.Cells.Item(145, 3).Value = 'The image data (A) includes a black pixel extending from the first pixel to the final pixel in the main scan direction.'
$ws.Cells.Item(146, 3).Value = 'Accordingly, all the light emitting elements 306 included in the light emitting unit 305 emit light.'
$ws.Cells.Item(147, 3).Value = 'An OLED head is an aspect of the light emitting unitFIG.10 is a diagram illustrating the cumulative light emission time of the light emitting element 306 after five copies of the image data (A) are printed.'
$ws.Cells.Item(148, 3).Value = 'In a case where a difference value between a maximum value of the cumulative light emission time and a minimum value of the cumulative light emission time is larger than the forced light emission threshold, the light emitting elements 306 in which the cumulative light emission time is shorter than the adjustment time is forced to emit light.'
$ws.Cells.Item(149, 3).Value = 'FIG.11 is a diagram illustrating a state where the cumulative light emission time of light emitting element 306 reaches the adjustment time after five copies of the image data (A) are printed.'
$ws.Cells.Item(150, 3).Value = 'The light emission control unit 314 makes the light emitting elements 306 in which the cumulative light emission time is shorter than the adjustment time emit light until the adjustment time.'
$ws.Cells.Item(151, 3).Value = 'Accordingly, all the light emitting elements 306 have the cumulative light emission time longer than or equal to the adjustment time.'
$ws.Cells.Item(152, 3).Value = 'FIG.12 is a diagram illustrating the cumulative light emission time of the light emitting element 306 in a case where image data (B) is printed.'
$ws.Cells.Item(153, 3).Value = 'Unlike a case where the image data (A) is printed, although the image data (B) is printed, there is the light emitting element 306 which does not emit light.'
$ws.Cells.Item(154, 3).Value = 'FIG.13 is a diagram illustrating the cumulative light emission time of the light emitting element 306 after 15 copies of the image data (B) are printed.'
$ws.Cells.Item(155, 3).Value = 'A light emitting element 306 shortly after five copies of the image data (A) are printed is used as the light emitting elementIn a case where a difference value between a maximum value of the cumulative light emission time and a minimum value of the cumulative light emission time is larger than the forced light emission threshold, the light emitting element 306 in which the cumulative light emission time is shorter than the adjustment time is forced to emit light.'
$ws.Cells.Item(156, 3).Value = 'FIG.14 is a diagram illustrating a state where the cumulative light emission time of the light emitting element 306 reaches the adjustment time after 15 copies of the image data (B) are printed.'
$ws.Cells.Item(157, 3).Value = 'The light emitting elements 306 in which the cumulative light emission time is shorter than the adjustment time emits light until the adjustment time.'
$ws.Cells.Item(158, 3).Value = 'Accordingly, all the light emitting elements 306 have the cumulative light emission time longer than or equal to the adjustment time.'
$ws.Cells.Item(159, 3).Value = 'As such, the light emission control unit 314 determines whether or not the value of the cumulative light emission time stored in the light emission time storage unit 301 is larger than the value of the adjustment time.'
$ws.Cells.Item(160, 3).Value = 'As a result of determination, in a case where the value of the cumulative light emission time is not larger than the value of the adjustment time, the light emission control unit 314 performs control such that the light emitting elements 306 emits light.'
$ws.Cells.Item(161, 3).Value = 'In contrast, in a case where the value of the cumulative light emission time is larger than the value of the adjustment time, the light emission control unit 314 performs control such that the light emitting elements 306 does not emit light.'
$ws.Cells.Item(162, 3).Value = 'Accordingly, the light emission control unit 314 makes the light emitting elements 306 in which the value of the cumulative light emission time is not larger than the value of the adjustment time emit light until the value of the adjustment time.'
$ws.Cells.Item(163, 3).Value = 'Hence, although the light emitting element 306 emits light until the largest value of the cumulative light emission time, the light emission time is shortened by a difference value between the value of the cumulative light emission time and the value of the adjustment time.'
$ws.Cells.Item(164, 3).Value = 'Accordingly, the light emitting elements 306 uniformly maintains the amount of emitted light and is prevented from being degraded due to light emission by the difference value.'
$ws.Cells.Item(165, 3).Value = 'The accompanying claims and their equivalents are intended to cover such forms or modifications as would fall within the scope and spirit of the inventions. WHAT IS CLAIMED IS:An image forming device comprising: a light emitting unit that includes a plurality of light emitting elements which form an electrostatic latent image on a photosensitive drum; an adjustment time storage unit that stores adjustment time which is shorter than cumulative light emission time of a light emitting element with longest cumulative light emission time; and a light emission control unit that performs control such that the light emitting element whose cumulative light emission time is shorter than the adjustment time which is stored in the adjustment time storage unit emits light until the adjustment time, in a case where a difference between cumulative light emission time of one of the plurality of light emitting elements and cumulative light emission time of another light emitting element satisfies a predetermined condition.The image forming device according to Claim 1, further comprising: a light emission time storage unit that stores the cumulative light emission times for each of the plurality of light emitting elements which are included in the light emitting unit.The image forming device according to Claim 1, wherein the one light emitting element and another light emitting element are respectively a light emitting element having a value of a longest cumulative light emission time and a light emitting element having a value of a shortest cumulative light emission time.The image forming device according to Claim 1, further comprising: a longest storage unit that stores a value of cumulative light emission time of the light emitting unit which satisfies a predetermined condition.The image forming device according to Claim 4, further comprising: a shortest storage unit that stores a value of cumulative light emission time of the light emitting unit which satisfies a predetermined condition.The image forming device according to Claim 5, further comprising: an adjustment time calculation unit that calculates the adjustment time.The image forming device according to Claim 6, wherein the adjustment time calculation unit calculates adjustment time on the basis of a difference between cumulative light emission time which is stored in the longest storage unit and cumulative light emission time which is stored in the shortest storage unit.The image forming device according to Claim 1, wherein the light emission control unit makes the light emitting element emit light after image formation processing is performed.The image forming device according to Claim 1, wherein the light emitting elements of the light emitting unit are organic ELs.An image forming method comprising: emitting light by including a plurality of light emitting elements which form an electrostatic latent image on a photosensitive drum; storing adjustment time which is shorter than cumulative light emission time of a light emitting element with longest cumulative light emission time; and controlling light emission such that the light emitting element whose cumulative light emission time is shorter than the adjustment time which is stored in the storing of the adjustment time emits light until the adjustment time, in a case where a difference between cumulative light emission time of one of the plurality of light emitting elements and cumulative light emission time of another light emitting element satisfies a predetermined condition.'
$ws.Cells.Item(166, 3).Value = 'ABSTRACT According to one embodiment, an image forming device includes a light emitting unit, an adjustment time storage unit, and light emission control unit.'
$ws.Cells.Item(167, 3).Value = 'The light emitting unit includes a plurality of light emitting elements which form an electrostatic latent image on a photosensitive drum.'
$ws.Cells.Item(168, 3).Value = 'The adjustment time storage unit stores adjustment time which is shorter than cumulative light emission time of a light emitting element with longest cumulative light emission time.'

# --- Append new row 169 (new sentence pair split out in the realignment) ---
$ws.Range("A168").Copy() | Out-Null
$ws.Range("A169").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(169, 1).Value = 168
$ws.Cells.Item(169, 2).Value = '発光制御部は、複数の発光素子のうち、一の発光素子の累積発光時間と他の発光素子の累積発光時間の差が所定の条件を満たす場合、調整時間記憶部に記憶する調整時間よりも累積発光時間が短い発光素子に対して調整時間まで発光素子を発光させるように制御する。'
$ws.Cells.Item(169, 3).Value = 'The light emission control unit performs control such that the light emitting element whose cumulative light emission time is shorter than the adjustment time which is stored in the adjustment time storage unit emits light until the adjustment time, in a case where a difference between cumulative light emission time of one of the plurality of light emitting elements and cumulative light emission time of another light emitting element satisfies a predetermined condition.'

Write-Host "Edit complete"